# Developer Guide diagram update:
#   "PersonListPanel" -> "Task" + "ListPanel"   (TaskListPanel)
#   "PersonCard"       -> "Task" + "Card"        (TaskCard)
#
# PowerPoint splits a run when only part of its text is retyped, so we
# replace just the leading "Person" characters of each label - this keeps
# the untouched remainder ("ListPanel"/"Card") as its own run, matching
# how the author actually edited the text in the UI class diagram.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if (-not $shp.HasTextFrame) {
        continue
    }

    $tr = $shp.TextFrame.TextRange
    $text = $tr.Text

    if ($text -eq "PersonListPanel" -or $text -eq "PersonCard") {
        # "Person" is always the first 6 characters of both labels.
        $prefix = $tr.Characters(1, 6)
        $prefix.Text = "Task"
    }
}
